# TestMatrix.xlsx update: rework the "AlexNet" test block into "Flickr",
# rename the old "Accuracy" header to "Crop Size", add the computed
# "mean image" crop-size column (G), and refresh the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2: F2 stays "Iteration"; G2 becomes "Crop Size" ---
$ws.Range("G2").Value = "Crop Size"

# --- GoogLeNet Cars block (rows 3-8): network name unchanged, but make
#     sure it reads "GoogLeNet Cars" (was already, kept for safety) ---
$ws.Range("C3").Value = "GoogLeNet Cars"

# --- AlexNet block (rows 9-14) becomes the Flickr block ---
$ws.Range("C9").Value = "Flickr"

# --- New column G: crop size per-row, derived from the input size (D) ---
$ws.Range("G3").Value = 496
$ws.Range("G4").Value = 496
$ws.Range("G5").Value = 248
$ws.Range("G6").Value = 248
$ws.Range("G7").Value = 124
$ws.Range("G8").Value = 124

$ws.Range("G9").Value = 496
$ws.Range("G10").Value = 496
$ws.Range("G11").Value = 248
$ws.Range("G12").Value = 248
$ws.Range("G13").Value = 124
$ws.Range("G14").Value = 124

# --- Selection now highlights the refreshed Flickr "Network" column ---
[void]$ws.Range("C9:C14").Select()
